# hackerrank_submissions_olaf_.xlsx - "Add files via upload" edit
#
# Changes applied:
#   1. Rename the worksheet from "Olaf's hackerrank submissions" to
#      "hackerrank submissions".
#   2. Header row (A1:F1) gets vertical-center + wrap-text alignment
#      (in addition to its existing medium-box border).
#   3. Column D width grows from 5 to 6 characters.
#   4. The "Challenge (link)" header text is expanded with a hint about
#      enabling editing in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet.
$ws.Name = "hackerrank submissions"

# 2. Header row formatting: vertical center + wrap text.
$headerRange = $ws.Range("A1:F1")
$headerRange.VerticalAlignment = -4108   # xlCenter
$headerRange.WrapText = $true

# 3. Widen column D (Points) from 5 to 6 characters.
$ws.Columns.Item(4).ColumnWidth = 5.166666666666667

# 4. Update the "Challenge (link)" header text.
$ws.Range("A1").Value = "Challenge (link - in Excel click 'enable editing' if not visible)"
